# Commit: "Include process sets splitting processes by sector"
#
# Adds 7 new TFM_Psets rows (one per TIMES-IE sector) to the "Sets-Proc"
# sheet, defining a process set ("All <SECTOR> processes") per sector
# using a wildcard PSET_PN criterion, gathered via AND/OR token settings
# matching the existing rows on the sheet. Also updates the workbook's
# active-sheet / selection bookkeeping (Sets-Proc becomes the active tab,
# SRV_Sets-Proc loses its tabSelected flag).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sets-Proc")

# --- New data rows 14-20 -------------------------------------------------
# Columns: A=PSET_SET  B=PSET_PN  F=SetName  G=SetDesc
#          H=T_Pos_AndOr  I=T_Neg_AndOr  J=T_Pos_AndOr_forSets  K=T_Neg_AndOr_forSets

$rows = @(
    @{ Row = 14; A = $null;   B = "A*,FT-AGR*";        F = "PRC_AGR"; G = "All AGR processes" },
    @{ Row = 15; A = $null;   B = "S-*,FT-SRV*";       F = "PRC_SRV"; G = "All SRV processes" },
    @{ Row = 16; A = "-IRE";  B = "I*,FT-IND*";        F = "PRC_IND"; G = "All IND processes" },
    @{ Row = 17; A = $null;   B = "P*,FT-PWR*,*GRID*"; F = "PRC_PWR"; G = "All PWR processes" },
    @{ Row = 18; A = $null;   B = "R*,FT-RSD*";        F = "PRC_RSD"; G = "All RSD processes" },
    @{ Row = 19; A = $null;   B = "S*,FT-SUP*,-S-*";   F = "PRC_SUP"; G = "All SUP processes" },
    @{ Row = 20; A = "-IRE";  B = "T*,FT-TRA*";        F = "PRC_TRA"; G = "All TRA processes" }
)

foreach ($r in $rows) {
    $n = $r.Row
    if ($r.A -ne $null) {
        # Leading apostrophe forces text/quote-prefix so a leading "-" is
        # not mistaken for a formula/negative number (matches the other
        # "-IRE" cell already on the sheet).
        $ws.Range("A$n").Value = "'" + $r.A
    }
    $ws.Range("B$n").Value = $r.B
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = "AND"
    $ws.Range("I$n").Value = "OR"
    $ws.Range("J$n").Value = "AND"
    $ws.Range("K$n").Value = "OR"
}

# --- Active sheet / selection bookkeeping --------------------------------

# SRV_Sets-Proc was the previously-active tab; move its cursor and let the
# tab-selected flag move to Sets-Proc below.
$wsSrv = $wb.Worksheets.Item("SRV_Sets-Proc")
[void]$wsSrv.Range("H31").Select()

[void]$ws.Activate()
[void]$ws.Range("B27").Select()
